# Upload new version with timestamp
# Adds 4 newly out-of-stock products (in their correct alphabetically-sorted
# position) to the "نواقص الأصناف" (missing items) report, recomputes the
# total, and refreshes the generated-at timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert 4 blank rows at the positions the new products belong to.
#    Working from the bottom of the table upward so earlier row numbers
#    used below are never invalidated by a prior insert.
# ---------------------------------------------------------------------
$ws.Rows("13:13").Insert()   # blank row that will hold SILVIRBURN
$ws.Rows("13:13").Insert()   # blank row that will hold TELFAST
$ws.Rows("12:12").Insert()   # blank row that will hold PANADOL
$ws.Rows("8:8").Insert()     # blank row that will hold DEPAKINE

# ---------------------------------------------------------------------
# 2. Re-apply the table's row formatting (styles + merged cells) to the
#    4 freshly inserted blank rows by copying it from row 7, which still
#    carries the canonical look of a product row.
# ---------------------------------------------------------------------
$templateRow = "7"
$newRows = @("8", "13", "15", "16")
$ws.Range("A" + $templateRow + ":Q" + $templateRow).Copy()
foreach ($r in $newRows) {
    $ws.Range("A" + $r + ":Q" + $r).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("A" + $r + ":B" + $r).Merge()
    $ws.Range("C" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
    $ws.Range("N" + $r + ":O" + $r).Merge()
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Write the 4 new products' data into place.
# ---------------------------------------------------------------------
$ws.Range("C8").Value = "DEPAKINE CHRONO 500MG 30 SCORED PROLONGED REL. F.C. TAB."
$ws.Range("H8").Value = "1:0"
$ws.Range("L8").Value = "1"
$ws.Range("N8").Value = "144.00"
$ws.Range("P8").Value = "144.0000"
$ws.Range("Q8").Value = "1:0"

$ws.Range("C13").Value = "PANADOL ACUTE HEAD COLD"
$ws.Range("H13").Value = "5:1"
$ws.Range("L13").Value = "0"
$ws.Range("N13").Value = "62.00"
$ws.Range("P13").Value = "31.0000"
$ws.Range("Q13").Value = "1:0"

$ws.Range("C15").Value = "SILVIRBURN CREAM 250 GM"
$ws.Range("H15").Value = "1:0"
$ws.Range("L15").Value = "1"
$ws.Range("N15").Value = "38.00"
$ws.Range("P15").Value = "38.0000"
$ws.Range("Q15").Value = "1:0"

$ws.Range("C16").Value = "TELFAST 30MG/5ML SUSP. 100 ML"
$ws.Range("H16").Value = "9:0"
$ws.Range("L16").Value = "1"
$ws.Range("N16").Value = "50.00"
$ws.Range("P16").Value = "50.0000"
$ws.Range("Q16").Value = "1:0"

# ---------------------------------------------------------------------
# 4. Renumber the "م" (row index) column so it stays a sequential
#    1..14 list spanning every product row (A7:A20).
# ---------------------------------------------------------------------
$firstRow = 7
$lastRow = 20
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("A" + $r).Value = ($r - $firstRow + 1)
}

# ---------------------------------------------------------------------
# 5. Update the grand-total cell (sum of the "سعر البيع" column), which
#    now sits on row 21 after the 4-row insertion.
# ---------------------------------------------------------------------
$ws.Range("P21").Value = 591.83

# ---------------------------------------------------------------------
# 6. Refresh the "generated at" timestamp in the footer (now row 22).
# ---------------------------------------------------------------------
$ws.Range("A22").Value = "Friday, 26 September, 2025 5:42 PM"
